# "Hours of Work.xlsx" — add a new "Debugging" sub-task (8 hours) to the
# "Expected Dev time" section, and bump the "Setup first machine" entry
# under "Test Time" from 0.25 to 3 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Expected Dev time" block (row 63) gets a new blank-row gap below it,
# just like every other section on this sheet, so insert 2 rows there.
# That pushes the old row 63 ("Expected Dev time" + its SUM) down to row 65,
# and everything below shifts down by 2 rows as well.
$ws.Rows("63:64").Insert()

# Row 63 is now free — fill it in with the new "Debugging" line item.
$ws.Range("A63").Value = "Debugging"
$ws.Range("C63").Value = 8

# "Expected Dev time" (now row 65) needs its SUM range extended so it picks
# up the newly-inserted row 63 (it used to be =SUM(C2:C62)).
$ws.Range("C65").Formula = "=SUM(C2:C64)"

# "Setup first machine" (now row 68, under the "EC2 / S3 - First machine
# setup" sub-task) increased from 0.25 to 3 hours.
$ws.Range("C68").Value = 3

# "Test Time" (now row 97) sums the shifted block - update its range too
# (it used to be =SUM(C65:C94)).
$ws.Range("C97").Formula = "=SUM(C67:C96)"

# "Total Time" (now row 101) references the other summary cells, which all
# moved down by 2 rows (it used to be =SUM(C97,C95,C63)).
$ws.Range("C101").Formula = "=SUM(C99,C97,C65)"

# Leave the view scrolled/selected roughly where the edit happened.
$ws.Range("C70").Select()
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
